$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove "sheet 2" (시트 2) entirely.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("시트 2")
$ws2.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. On the remaining sheet, insert a new column D ("comment"), shifting the
#    old D (float_val/float) and E (bool_val/bool) columns one place right
#    (-> E and F).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("sheet 1")
$ws.Columns.Item(4).Insert()

# The insert operation re-computed the two originally-empty cells (old D6 and
# old E7, now at E6/F7) as 0 instead of leaving them blank - put them back.
$ws.Range("E6").ClearContents()
$ws.Range("F7").ClearContents()

# ---------------------------------------------------------------------------
# 3. Populate the new column D with its header/type row + data.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "안녕"
$ws.Range("D2").Value = "comment"
$ws.Range("D3").Value = 111
$ws.Range("D4").Value = "攻击"
$ws.Range("D5").Value = 2345
$ws.Range("D6").Value = 9999
$ws.Range("D7").Value = 1111

# ---------------------------------------------------------------------------
# 4. Re-label the two header cells whose text changed.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "어머나"
$ws.Range("F1").Value = "闪避"
